# Apply cell value updates from the crypto price refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '26.662.67'
    'E2' = '  +1.22%  '
    'D3' = '1.631.08'
    'E3' = '  +0.57%  '
    'E4' = '  +0.02%  '
    'D5' = '213.36'
    'E5' = '  +0.59%  '
    'E6' = '  +3.30%  '
    'E7' = '  +0.03%  '
    'E8' = '  +1.81%  '
    'E9' = '  +0.93%  '
    'D10' = '19.21'
    'E10' = '  +2.89%  '
    'E11' = '  +3.18%  '
    'D12' = '1.858.01'
    'E12' = '  +0.48%  '
    'D13' = '1.609.87'
    'E13' = '  -0.79%  '
    'E14' = '  +1.44%  '
    'E15' = '  +1.12%  '
    'D16' = '26.654.55'
    'E16' = '  +1.14%  '
    'D17' = '63.48'
    'E17' = '  +1.72%  '
    'D18' = '0.0₃0740'
    'E18' = '  +2.15%  '
    'D19' = '218.73'
    'E19' = '  +8.06%  '
    'E20' = '  +0.03%  '
    'D21' = '4.29'
    'E21' = '  +0.76%  '
    'E22' = '  +2.00%  '
    'D24' = '1.95'
    'E24' = '  +4.57%  '
    'D25' = '147.81'
    'E25' = '  +2.16%  '
    'E26' = '  +0.00%  '
    'E27' = '  +1.20%  '
    'E28' = '  +4.23%  '
    'D29' = '15.53'
    'E29' = '  +2.25%  '
    'E30' = '  -3.07%  '
    'E31' = '  +0.34%  '
    'E32' = '  +3.84%  '
    'E33' = '  +2.20%  '
    'E34' = '  +0.97%  '
    'E35' = '  +0.35%  '
    'D36' = '1.214.40'
    'E36' = '  +4.97%  '
    'E37' = '  +4.44%  '
    'D38' = '0.805'
    'E38' = '  +0.13%  '
    'E39' = '  +0.06%  '
    'D40' = '0.500'
    'E40' = '  +0.48%  '
    'E41' = '  -1.81%  '
    'D42' = '0.795'
    'E42' = '  +1.50%  '
    'D43' = '5.34'
    'E43' = '  -0.75%  '
    'D44' = '1.768.70'
    'E44' = '  +0.45%  '
    'D45' = '92.77'
    'E45' = '  +0.12%  '
    'E46' = '  +2.58%  '
    'B47' = 'BabyDogeCoin'
    'C47' = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
    'D47' = '0.0₆0105'
    'E47' = '  +0.17%  '
    'B48' = 'Aave'
    'C48' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    'D48' = '55.09'
    'E48' = '  +2.45%  '
    'E49' = '  +0.48%  '
    'D50' = '7.61'
    'E50' = '  +4.25%  '
    'E51' = '  -0.06%  '
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    $cell.Style = "Normal"
}
